# Strip the thousands-separator comma from the PATH-ID labels stored in
# column A of the PATH sheet (e.g. "1,001" -> "1001"). The column holds
# text labels (not numbers), so the cells must stay text after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PATH")

$lastRow = 41

# Force the column to remain text before writing back so Excel doesn't
# reinterpret the comma-free digits as a real number.
$ws.Range("A1:A41").NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($current -ne $null -and $current.Contains(",")) {
        $cell.Value = $current.Replace(",", "")
    }
}
